$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1.247512323997491
$ws.Range("E2").Value = 2.82953744009995
$ws.Range("C3").Value = 0.5917823527752164
$ws.Range("E3").Value = -1.985049937499994
$ws.Range("C4").Value = 1.985742476411234
$ws.Range("E4").Value = 7.819356632099961
$ws.Range("C5").Value = -0.7518797681958955
$ws.Range("E5").Value = -1.590425574400001
$ws.Range("C6").Value = -0.5765930039052902
$ws.Range("E6").Value = -1.590425574400001
$ws.Range("C7").Value = -0.07642926654481963
$ws.Range("E7").Value = -0.3994003998999851
$ws.Range("C8").Value = 0.9274109147535681
$ws.Range("E8").Value = 2.829537440099972
$ws.Range("C9").Value = -0.5259734324162268
$ws.Range("E9").Value = -1.194610791900008
$ws.Range("C10").Value = 0.07456754038981384
$ws.Range("E10").Value = 0.4006004000999708
$ws.Range("C11").Value = -0.07666472728168339
$ws.Range("E11").Value = -1.590425574400012
$ws.Range("C12").Value = -0.2004754673795017
$ws.Range("E12").Value = -0.3994003998999962
$ws.Range("C13").Value = -0.7283174404322912
$ws.Range("E13").Value = -2.378486270400004
$ws.Range("C14").Value = -0.03096525636255842
$ws.Range("E14").Value = 1.205410808099971
$ws.Range("C15").Value = -0.04074803603358879
$ws.Range("E15").Value = -0.2354831990173722
$ws.Range("C16").Value = 1.296559640836992
$ws.Range("E16").Value = 0.439921111559638
$ws.Range("C17").Value = 0.8143067496459322
$ws.Range("E17").Value = -0.05551881767973388
$ws.Range("C18").Value = -1.077932696718564
$ws.Range("E18").Value = -0.02296885644381685
$ws.Range("C19").Value = 0.4743170952486997
$ws.Range("E19").Value = -0.324019587165425
